$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily-push row (row 97) with the same layout as the
# existing data rows: date, weekday, hour, ranking.
$ws.Range("A97").Value = "'2025/10/13"
$ws.Range("A97").Style = "Normal"
$ws.Range("B97").Value = "月"
$ws.Range("C97").Value = 5
$ws.Range("D97").Value = 49
